# Commit: "Change Excel Field View to Cache, And set default value to FALSE"
#
# The workbook's per-sheet "Property" table has a shared column historically
# labeled "View" (a TRUE/FALSE flag). This rewrites that header to "Cache"
# on every sheet that has it, and resets every existing TRUE value in that
# column back to the new default of FALSE.

$wb = $excel.ActiveWorkbook

$activeSheetName = "Record_BagItemList"

foreach ($ws in $wb.Worksheets) {

    $header = $ws.Rows.Item(1).Find("View")
    if ($header -eq $null) {
        continue
    }

    $col = $header.Column
    $lastRow = $ws.UsedRange.Rows.Count

    # Flip every TRUE in the column back to the new FALSE default.
    if ($lastRow -ge 2) {
        $dataRange = $ws.Cells.Item(2, $col).Resize($lastRow - 1, 1)

        for ($r = 2; $r -le $lastRow; $r++) {
            $cell = $ws.Cells.Item($r, $col)
            if ($cell.Value2 -eq $true) {
                $cell.Value = $false
            }
        }

        # The Property sheet's column was regenerated wholesale by the
        # data-export tool, which dropped the few per-row highlight styles
        # that had been manually applied on top of the boolean cells.
        if ($ws.Name -eq "Property") {
            [void]$dataRange.ClearFormats()
        }
    }

    # Rename the column header itself.
    $header.Value = "Cache"

    [void]$header.Select()
}

# The data-export run ended with this sheet on screen.
$active = $wb.Worksheets.Item($activeSheetName)
$active.Activate()
[void]$active.Range("G1").Select()
